$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.571.60"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.598.65"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  -0.10%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "609.36"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "148.79"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -1.27%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "8.07"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("E10").Value = "  -0.13%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.416"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "4.204.19"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  +0.99%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "29.85"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "3.559.04"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "66.670.94"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("E17").Value = "  +0.86%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "11.51"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.35%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.38"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.95%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "15.12"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "428.44"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.61%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.619"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "78.87"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "3.746.25"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +3.89%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.34"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.29%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.56"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "3.593.65"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("E33").Value = "  +2.39%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "25.46"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.22%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "7.87"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("E38").Value = "  -1.91%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "177.07"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.11%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0858"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.54%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.97%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.899"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("E44").Value = "  +9.11%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").Value = "  -1.94%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "25.01"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.01%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "24.03"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.27%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.19"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.83%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.954"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "2.423.66"
$ws.Range("E51").Value = "  +5.09%  "
